$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Riteriai", 3, 0, 4, 4, 7, 4, 4, 6, 8, 6, 4, 5),
    @("Šiauliai", 4, 3, 4, 2, 4, 5, 2, 3, 1, 9, 1, 4),
    @("Hegelmann Litauen", 2, 2, 4, 8, 5, 5, 5, 2, 6, 3, 3, 4),
    @("Kauno Žalgiris", 1, 4, 5, 8, 5, 5, 0, 0, 2, 1, 4, 2),
    @("Džiugas Telšiai", 1, 1, 4, 1, 4, 2, 2, 2, 5, 0, 3, 1),
    @("Banga", 1, 1, 4, 1, 1, 7, 1, 1, 4, 2, 4, 6),
    @("Sūduva", 4, 4, 8, 3, 0, 4, 1, 2, 1, 1, 6, 2),
    @("Žalgiris", 2, 7, 2, 2, 4, 2, 3, 3, 6, 4, 1, 4),
    @("Panevėžys", 1, 4, 3, 5, 3, 7, 3, 2, 3, 3, 3, 5),
    @("Dainava", 3, 0, 4, 1, 1, 3, 1, 5, 6, 6, 5, 11)
)

$rowIndex = 2
foreach ($rowValues in $data) {
    $colIndex = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
